$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, shifting existing rows 111..239 down to 112..240.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly price record.
$ws.Range("A111").Value = 5
$ws.Range("B111").Value = "Macroferia Regional de Talca"
$ws.Range("C111").Value = "Maule"
$ws.Range("D111").Value = 44629
$ws.Range("E111").Value = 7
$ws.Range("F111").Value = 100112006
$ws.Range("G111").Value = "Repollo"
$ws.Range("H111").Value = "Crespo record"
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 1000
$ws.Range("L111").Value = 1000
$ws.Range("M111").Value = 1000
$ws.Range("N111").Value = "$/unidad"
$ws.Range("O111").Value = "Región del Maule"
$ws.Range("P111").Value = 1000
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = "Hortaliza"
